# Update for "Add data for 2022-11-22" (carjacking-by-neighborhood-by-month)
# - Rename the sheet / update the "through" label from Nov 13 to Nov 14
# - Bump a handful of neighborhood/month carjacking counts

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet tab name
$ws.Name = "Through 2022-11-14"

# Header cell describing the cutoff date for the current month's column
$ws.Range("B1").Value = "November 2022 (through November 14)"

# Row 2 - Englewood
$ws.Range("M2").Value = 5

# Row 3 - Auburn Gresham
$ws.Range("CA3").Value = 4

# Row 5 - Garfield Park
$ws.Range("B5").Value = 5
$ws.Range("X5").Value = 7
$ws.Range("AT5").Value = 3
$ws.Range("BE5").Value = 3

# Row 6 - Humboldt Park
$ws.Range("B6").Value = 2
$ws.Range("M6").Value = 4
$ws.Range("BP6").Value = 4

# Row 9 - Grand Crossing
$ws.Range("M9").Value = 3

# Row 14 - Kenwood
$ws.Range("X14").Value = 1
$ws.Range("BE14").Value = 1

# Row 16 - North Lawndale
$ws.Range("B16").Value = 1
$ws.Range("CA16").Value = 2

# Row 20 - West Pullman
$ws.Range("X20").Value = 2

# Row 25 - Austin
$ws.Range("M25").Value = 6

# Row 26 - Bridgeport
$ws.Range("M26").Value = 3

# Row 36 - Belmont Cragin
$ws.Range("BP36").Value = 1

# Row 38 - Irving Park
$ws.Range("B38").Value = 2

# Row 40 - Roseland
$ws.Range("B40").Value = 3
$ws.Range("BP40").Value = 2

# Row 54 - Beverly
$ws.Range("X54").Value = 1
$ws.Range("BE54").Value = 1

# Row 64 - Gage Park
$ws.Range("B64").Value = 1

# Row 75 - Lincoln Park
$ws.Range("X75").Value = 4

# Row 76 - Little Italy, UIC
$ws.Range("M76").Value = 2

# Row 81 - Morgan Park
$ws.Range("AI81").Value = 2

# Row 97 - Ukrainian Village
$ws.Range("X97").Value = 1
